# Refresh the cryptos price list (GitHub Actions scheduled update).
# Note: several "Price" (column D) values look numeric (e.g. 53.32, 0.0740,
# 1.00) but must stay stored as literal text, matching the source data
# (trailing zeros / exact formatting must be preserved, e.g. "0.0740" not
# "0.074", "1.00" not "1"). A leading "'" forces Excel to keep them as text
# instead of silently coercing to a Double and losing formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.406.54"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").Value = "2.040.93"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'245.19"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'53.32"
$ws.Range("E8").Value = "  -7.48%  "
$ws.Range("D9").Value = "'61.98"
$ws.Range("E9").Value = "  +4.64%  "
$ws.Range("E10").Value = "  -3.29%  "
$ws.Range("D11").Value = "'0.0740"
$ws.Range("E11").Value = "  -5.03%  "
$ws.Range("E12").Value = "  -4.14%  "
$ws.Range("D13").Value = "'0.923"
$ws.Range("E13").Value = "  +5.81%  "
$ws.Range("D14").Value = "'14.40"
$ws.Range("E14").Value = "  -5.37%  "
$ws.Range("D15").Value = "2.340.81"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").Value = "'5.35"
$ws.Range("E16").Value = "  -4.86%  "
$ws.Range("D17").Value = "2.035.81"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "36.360.47"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("D19").Value = "'16.88"
$ws.Range("E19").Value = "  -5.62%  "
$ws.Range("D20").Value = "'71.07"
$ws.Range("E20").Value = "  -3.43%  "
$ws.Range("D21").Value = "0.0₃0848"
$ws.Range("E21").Value = "  -4.72%  "
$ws.Range("D22").Value = "'235.46"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'5.13"
$ws.Range("E23").Value = "  -4.91%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "'2.36"
$ws.Range("E25").Value = "  -3.47%  "
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("D27").Value = "'163.56"
$ws.Range("E27").Value = "  -2.96%  "
$ws.Range("D28").Value = "'9.04"
$ws.Range("E28").Value = "  -11.99%  "
$ws.Range("D29").Value = "'19.72"
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("E30").Value = "  -2.67%  "
$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = "  +5.44%  "
$ws.Range("D32").Value = "'4.96"
$ws.Range("E32").Value = "  -9.30%  "
$ws.Range("E33").Value = "  -4.42%  "
$ws.Range("D34").Value = "'4.35"
$ws.Range("E34").Value = "  -7.38%  "
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.0864"
$ws.Range("E36").Value = "  +5.38%  "
$ws.Range("D37").Value = "'1.82"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("E38").Value = "  -6.22%  "
$ws.Range("D39").Value = "'4.93"
$ws.Range("E39").Value = "  -4.17%  "
$ws.Range("E40").Value = "  -7.46%  "
$ws.Range("E41").Value = "  -4.42%  "
$ws.Range("D42").Value = "'0.0211"
$ws.Range("E42").Value = "  -5.57%  "
$ws.Range("E43").Value = "  -4.94%  "
$ws.Range("D44").Value = "'92.61"
$ws.Range("E44").Value = "  -4.17%  "
$ws.Range("D45").Value = "'0.0886"
$ws.Range("D46").Value = "1.376.42"
$ws.Range("E46").Value = "  +5.71%  "
$ws.Range("D47").Value = "'15.52"
$ws.Range("E47").Value = "  -8.17%  "
$ws.Range("D48").Value = "'7.36"
$ws.Range("E48").Value = "  +9.11%  "
$ws.Range("E49").Value = "  +1.42%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.224.32"
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'2.23"
$ws.Range("E51").Value = "  -5.55%  "
